$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
$newDate = Get-Date -Year 2023 -Month 9 -Day 8 -Hour 0 -Minute 0 -Second 0
$ws.Range("C2:C$lastRow").Value = $newDate.Date
